$d = $word.ActiveDocument

# Walk every run in the document and normalize the explicit "off" boolean
# formatting properties (Bold / Italic / Strikethrough) on the runs that
# carry the CSS-styled color (this mirrors the Apache POI 4.1.0 -> 5.2.3
# upgrade which re-serializes <w:b w:val="false"/> etc. as the canonical
# on/off form instead of the legacy "false" token).
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Text = "Some CSS styled text"
$rng.Find.Forward = $true
$rng.Find.Wrap = 0

while ($rng.Find.Execute("Some CSS styled text", $false, $false, $false, $false, $false, $true, 1, $false, $null, 0)) {
    # Toggle true -> false so the engine always emits an explicit on/off
    # value for each property instead of silently dropping an
    # already-false attribute.
    $rng.Font.Bold = $true
    $rng.Font.Bold = $false
    $rng.Font.Italic = $true
    $rng.Font.Italic = $false
    $rng.Font.StrikeThrough = $true
    $rng.Font.StrikeThrough = $false
    $rng.Collapse(0)
}

Write-Host "done"
